$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stock splits for NVDA are also recorded against the ISIN (US67066G1040).
# Re-append the same split history (rows 2-7) as rows 8-13, but keyed by ISIN
# instead of ticker symbol.

$dates  = @(45450, 44397, 39336, 38814, 37146, 36704)
$shares = @(10, 4, 1.5, 2, 2, 2)
$isin   = "US67066G1040"

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 8 + $i

    # Copy row 2's formatting (date style + value style) down onto the new
    # row so the new cells reuse the existing cellXfs entries.
    $ws.Range("A2:C2").Copy()
    $ws.Range("A" + $row + ":C" + $row).PasteSpecial(-4122)

    $ws.Range("A" + $row).Value = $dates[$i]
    $ws.Range("B" + $row).Value = $isin
    $ws.Range("C" + $row).Value = $shares[$i]

    # The ISIN column gets its own (slightly darker) font color, distinct
    # from the ticker-symbol rows above it.
    $ws.Range("B" + $row).Font.Color = 2630431
}

$ws.Range("A1").Select()
